# Daily attendance processing - 2026-01-07 06:45:59
# Normalize the "Recorded By" (column G) values by reversing the order of the
# comma-separated contributor list (e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val -notmatch ",") { continue }

    $parts = $val -split ", "
    if ($parts.Count -lt 2) { continue }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.ToLower() -eq "system") { $hasSystem = $true }
    }
    if (-not $hasSystem) { continue }

    $reversed = @()
    for ($i = $parts.Count - 1; $i -ge 0; $i--) {
        $reversed += $parts[$i]
    }
    $newVal = [string]::Join(", ", $reversed)

    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}
